$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 363, pushing existing rows 363-366 down to 365-368
$ws.Rows.Item(363).Resize(2).Insert()

# Copy the style of column D (date column) down to the new rows so the "s=2" style carries over
$ws.Range("D363").Value = 44628
$ws.Range("D363").NumberFormat = $ws.Range("D365").NumberFormat

$ws.Range("D364").Value = 44628
$ws.Range("D364").NumberFormat = $ws.Range("D365").NumberFormat

# Row 363: new weekly data
$ws.Range("A363").Value = 11
$ws.Range("B363").Value = "Vega Monumental Concepción"
$ws.Range("C363").Value = "Bíobío"
$ws.Range("E363").Value = 8
$ws.Range("F363").Value = 100112020
$ws.Range("G363").Value = "Tomate"
$ws.Range("H363").Value = "Larga vida"
$ws.Range("I363").Value = "Primera"
$ws.Range("J363").Value = 220
$ws.Range("K363").Value = 17000
$ws.Range("L363").Value = 18000
$ws.Range("M363").Value = 17545
$ws.Range("N363").Value = "$/bandeja 18 kilos"
$ws.Range("O363").Value = "Región del Maule"
$ws.Range("P363").Value = 975
$ws.Range("Q363").Value = 18
$ws.Range("R363").Value = "Hortaliza"

# Row 364: new weekly data
$ws.Range("A364").Value = 11
$ws.Range("B364").Value = "Vega Monumental Concepción"
$ws.Range("C364").Value = "Bíobío"
$ws.Range("E364").Value = 8
$ws.Range("F364").Value = 100112020
$ws.Range("G364").Value = "Tomate"
$ws.Range("H364").Value = "Semiduro"
$ws.Range("I364").Value = "Primera"
$ws.Range("J364").Value = 250
$ws.Range("K364").Value = 5000
$ws.Range("L364").Value = 5500
$ws.Range("M364").Value = 5300
$ws.Range("N364").Value = "$/caja 10 kilos"
$ws.Range("O364").Value = "Quillón"
$ws.Range("P364").Value = 530
$ws.Range("Q364").Value = 10
$ws.Range("R364").Value = "Hortaliza"

$wb.Save()
